$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A (was 9.9 chars -> now ~13.2 chars)
$ws.Columns.Item(1).ColumnWidth = 12.3

# Materialize empty rows 6-9 (no data, just present in sheetData)
$ws.Rows.Item(6).Hidden = $true
$ws.Rows.Item(6).Hidden = $false
$ws.Rows.Item(7).Hidden = $true
$ws.Rows.Item(7).Hidden = $false
$ws.Rows.Item(8).Hidden = $true
$ws.Rows.Item(8).Hidden = $false
$ws.Rows.Item(9).Hidden = $true
$ws.Rows.Item(9).Hidden = $false

# New block of data in rows 10-12
$ws.Range("A10").Value = "ACDC"
$ws.Range("B10").Value = "BTO"

$ws.Range("A11").Value = "Align Tech"
$ws.Range("B11").Value = "Boeing"
$ws.Range("C11").Value = "Citigroup"

$ws.Range("A12").Value = "Audi"
$ws.Range("B12").Value = "Buick"
